# Fruta / hortaliza, semanal
# Insert a new weekly record at row 57 of the "Coco" (Vega Modelo de Temuco)
# price sheet, pushing all subsequent rows (57-121) down by one (58-122).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 57; this shifts rows 57:121
# down to 58:122, preserving all of their data and formatting.
$ws.Rows(57).Insert()

# Populate the newly inserted row 57 with this week's record.
$ws.Range("A57").Value = 10
$ws.Range("B57").Value = "Vega Modelo de Temuco"
$ws.Range("C57").Value = "La Araucanía"
$ws.Range("D57").Value = 45159
$ws.Range("E57").Value = 9
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100108
$ws.Range("H57").Value = "Tropicales y subtropicales"
$ws.Range("I57").Value = 100108007
$ws.Range("J57").Value = "Coco"
$ws.Range("K57").Value = "Sin especificar"
$ws.Range("L57").Value = "Primera"
$ws.Range("M57").Value = 50
$ws.Range("N57").Value = 36000
$ws.Range("O57").Value = 36000
$ws.Range("P57").Value = 36000
$ws.Range("Q57").Value = "$/malla 20 unidades"
$ws.Range("R57").Value = "Perú"
$ws.Range("S57").Value = 1800
$ws.Range("T57").Value = 20
